# queuing_theory.xlsx -- re-run the five Erlang/queueing worksheets with new
# inputs ("change EN calculation, split erlang c into two programs").
#
# Every input/result cell in this workbook is stored as literal TEXT (not a
# number) even though the text looks numeric. Writing through Range.Value
# directly would let Excel's normal type-inference turn "0.005" into the
# number 0.005, which would not match the source data's cell typing. To
# keep the cells as plain text (no helper styles / quote-prefix markers
# either), we stage each new value as a text formula in a scratch cell,
# copy it, and Paste-Special "Values only" onto the destination -- this
# pastes the literal text without carrying the formula or any formatting.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Sheet,
        [string]$Address,
        [string]$Text
    )
    $scratch = $Sheet.Range("ZZ1")
    # Quote-as-text formula so the staged value is TEXT, not a number.
    $escaped = $Text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $Sheet.Range($Address).PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# ---------------------------------------------------------------- Part 1 --
$ws1 = $wb.Worksheets.Item(1)
Set-TextValue $ws1 "B2" "0.005"
Set-TextValue $ws1 "D2" "15"
Set-TextValue $ws1 "B3" "1.75"
Set-TextValue $ws1 "D3" "0.0033186086434399764"
Set-TextValue $ws1 "B4" "0.25"

# ---------------------------------------------------------------- Part 2 --
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").ClearContents()
Set-TextValue $ws2 "D3" "0.1018867924528302"
Set-TextValue $ws2 "B4" "15.0"
Set-TextValue $ws2 "D4" "4.528301886792453"

# ---------------------------------------------------------------- Part 3 --
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3 "B2" "30.0"
Set-TextValue $ws3 "D2" "0.8196836008567772"
Set-TextValue $ws3 "B3" "15.0"
Set-TextValue $ws3 "D3" "81.14237961445023"
Set-TextValue $ws3 "B4" "0.5"
Set-TextValue $ws3 "D4" "0.012324016189216061"
Set-TextValue $ws3 "D5" "0.9753519676215678"

# ---------------------------------------------------------------- Part 4 --
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4 "B2" "4.0"
Set-TextValue $ws4 "D2" "0.6893203883495146"
Set-TextValue $ws4 "B3" "20.0"
Set-TextValue $ws4 "D3" "1.2427184466019416"
Set-TextValue $ws4 "B4" "20.0"

# ---------------------------------------------------------------- Part 5 --
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5 "B2" "1.0"
Set-TextValue $ws5 "D2" "0.5294117647058824"
Set-TextValue $ws5 "B3" "2.0"
Set-TextValue $ws5 "D3" "0.47058823529411764"
Set-TextValue $ws5 "B4" "1.0"
